$d = $word.ActiveDocument

foreach ($story in $d.StoryRanges) {
    $r = $story
    while ($r -ne $null) {
        $r.Find.Execute("вшэ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "тестик", 2) | Out-Null
        $r.Find.Execute("факультет", $true, $false, $false, $false, $false,
                         $true, 1, $false, "тестик", 2) | Out-Null
        $r.Find.Execute("/raise_error", $true, $false, $false, $false, $false,
                         $true, 1, $false, "тестик", 2) | Out-Null
        $r = $r.NextStoryRange
    }
}
